$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '37.397.12'
$ws.Range('E2').Value = '  +5.43%  '

# Row 3
$ws.Range('D3').Value = '2.040.28'
$ws.Range('E3').Value = '  +2.95%  '

# Row 4
$ws.Range('E4').Value = '  -0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.76%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.649'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.78%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '65.06'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +17.50%  '

# Row 8
$ws.Range('E8').Value = '  -0.03%  '

# Row 9
$ws.Range('E9').Value = '  +6.31%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '59.02'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.21%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0755'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.70%  '

# Row 12
$ws.Range('E12').Value = '  +0.90%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.906'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.83%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.10'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.31%  '

# Row 15
$ws.Range('D15').Value = '2.339.20'
$ws.Range('E15').Value = '  +3.03%  '

# Row 16
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.61'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +8.33%  '

# Row 17
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '20.83'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +23.19%  '

# Row 18
$ws.Range('D18').Value = '2.042.93'
$ws.Range('E18').Value = '  +3.05%  '

# Row 19
$ws.Range('D19').Value = '37.313.41'
$ws.Range('E19').Value = '  +5.30%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '73.11'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.08%  '

# Row 21
$ws.Range('D21').Value = '0.0₃0874'
$ws.Range('E21').Value = '  +5.45%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.35'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +8.05%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.61'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.60%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.76'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +23.32%  '

# Row 25
$ws.Range('E25').Value = '  -0.20%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.35'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.16%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.59'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +6.31%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '165.89'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.02%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.86'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.94%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.121'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.59%  '

# Row 31
$ws.Range('E31').Value = '  +9.73%  '

# Row 32
$ws.Range('E32').Value = '  +9.02%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.113'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +27.24%  '

# Row 34
$ws.Range('E34').Value = '  +11.74%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0616'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.47%  '

# Row 36
$ws.Range('E36').Value = '  +13.10%  '

# Row 37
$ws.Range('E37').Value = '  +0.06%  '

# Row 38
$ws.Range('E38').Value = '  +1.12%  '

# Row 39
$ws.Range('E39').Value = '  +23.53%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.103'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +18.91%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.23'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.84%  '

# Row 42
$ws.Range('E42').Value = '  +4.24%  '

# Row 43
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0219'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.80%  '

# Row 44
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.73'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +22.18%  '

# Row 45
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.14'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +11.29%  '

# Row 46
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.14'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.59%  '

# Row 47
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.11'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +12.19%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '95.48'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.12%  '

# Row 49
$ws.Range('D49').Value = '1.415.53'
$ws.Range('E49').Value = '  +2.83%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.93'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.69%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '47.50'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.35%  '
